$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: create a brand-new cell with the correct column style (1/2/3) by
# copying formats from a stable template cell in the same column, then
# writing the value. This avoids bloating styles.xml with new cellXfs.
# ---------------------------------------------------------------------------
function Set-NewCell($row, $col, $templateRow, $value) {
    $ws.Cells.Item($templateRow, $col).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, $col).Value = $value
}

# Column letters: 1=A (style 1, bold label), 2=B (style 2), 3=C (style 3, red)
# Row 10 (Objetivos:/...) keeps all three columns populated throughout, so it
# is used as the stable template row for new cells.

# --- Row 10: professor name replaces the old "Objetivos" long paragraph ----
$ws.Cells.Item(10,2).Value = "5840963 - Daniela Camargo Vernilli"
$ws.Cells.Item(10,3).Value = "5840963 - Daniela Camargo Vernilli"

# --- Row 13: new A13 "Programa resumido:" label; B13/C13 become "Semestral" -
$ws.Rows.Item(13).RowHeight = 60
Set-NewCell 13 1 10 "Programa resumido:"
$ws.Cells.Item(13,2).Value = "Semestral"
$ws.Cells.Item(13,3).Value = "Semestral"

# --- Row 14: A14 becomes "Short syllabus:"; B14/C14 removed ---------------
$ws.Rows.Item(14).RowHeight = 60
$ws.Cells.Item(14,1).Value = "Short syllabus:"
$ws.Cells.Item(14,2).Clear() | Out-Null
$ws.Cells.Item(14,3).Clear() | Out-Null

# --- Row 15: A15 becomes "Programa:"; B15/C15 created with "01/01/2012" ----
$ws.Rows.Item(15).RowHeight = 120
$ws.Cells.Item(15,1).Value = "Programa:"
Set-NewCell 15 2 10 "01/01/2012"
Set-NewCell 15 3 10 "01/01/2012"

# --- Row 16: A16 becomes "Syllabus:"; B16/C16 removed ----------------------
$ws.Rows.Item(16).RowHeight = 120
$ws.Cells.Item(16,1).Value = "Syllabus:"
$ws.Cells.Item(16,2).Clear() | Out-Null
$ws.Cells.Item(16,3).Clear() | Out-Null

# --- Row 17: A17 becomes "Avaliação:" --------------------------------------
$ws.Rows.Item(17).RowHeight = 0
$ws.Rows.Item(17).Height = $ws.Rows.Item(12).Height
$ws.Cells.Item(17,1).Value = "Avaliação:"

# --- Row 18: A18 becomes "Método:"; B18/C18 created with the professor text -
$ws.Rows.Item(18).RowHeight = 60
$ws.Cells.Item(18,1).Value = "Método:"
Set-NewCell 18 2 10 "5840963 - Daniela Camargo Vernilli"
Set-NewCell 18 3 10 "5840963 - Daniela Camargo Vernilli"

# --- Row 19: A19 "Critério:"; B19/C19 "Aulas expositivas..." --------------
$ws.Cells.Item(19,1).Value = "Critério:"
$ws.Cells.Item(19,2).Value = "Aulas expositivas, demonstrações, aulas de laboratório e projetos."
$ws.Cells.Item(19,3).Value = "Aulas expositivas, demonstrações, aulas de laboratório e projetos."

# --- Row 20: A20 "Norma de recuperação:"; B20/C20 "Média ponderada..." ----
$ws.Cells.Item(20,1).Value = "Norma de recuperação:"
$ws.Cells.Item(20,2).Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Cells.Item(20,3).Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

# --- Row 21: A21 "Bibliografia:"; B21/C21 "Aplicação de uma prova..." -----
$ws.Rows.Item(21).RowHeight = 120
$ws.Cells.Item(21,1).Value = "Bibliografia:"
$ws.Cells.Item(21,2).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Cells.Item(21,3).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# --- Row 22: A22 becomes "Requisitos:"; B22/C22 removed --------------------
$ws.Rows.Item(22).RowHeight = 0
$ws.Rows.Item(22).Height = $ws.Rows.Item(12).Height
$ws.Cells.Item(22,1).Value = "Requisitos:"
$ws.Cells.Item(22,2).Clear() | Out-Null
$ws.Cells.Item(22,3).Clear() | Out-Null

# --- Row 23: A23 removed; B23/C23 created with "LOQ4031..." ---------------
$ws.Rows.Item(23).RowHeight = 30
$ws.Cells.Item(23,1).Clear() | Out-Null
Set-NewCell 23 2 10 "LOQ4031 -  Química Geral I  (Requisito)`n"
Set-NewCell 23 3 10 "LOQ4031 -  Química Geral I  (Requisito)`n"

# --- Row 24: entirely removed (sheet now ends at row 23) ------------------
$ws.Rows.Item(24).Delete() | Out-Null

Write-Output "edit complete"
